$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 8): bold "title"-style labels in B8:D8 ---
$ws.Cells.Item(8, 2).Value = "Number of employees"
$ws.Cells.Item(8, 2).Font.Bold = $true

$ws.Cells.Item(8, 3).Value = "Assets (local currency, unless noted otherwise)"
$ws.Cells.Item(8, 3).Font.Bold = $true

$ws.Cells.Item(8, 4).Value = "Turnover (local currency, unless noted otherwise)"
$ws.Cells.Item(8, 4).Font.Bold = $true

# --- Row 9: Micro ---
$ws.Cells.Item(9, 1).Value = "Micro"
$ws.Cells.Item(9, 2).Value = "'"
$ws.Cells.Item(9, 2).Style = "Normal"
$ws.Cells.Item(9, 3).Value = "'"
$ws.Cells.Item(9, 3).Style = "Normal"
$ws.Cells.Item(9, 4).Value = "'"
$ws.Cells.Item(9, 4).Style = "Normal"

# --- Row 10: Small ---
$ws.Cells.Item(10, 1).Value = "Small"
$ws.Cells.Item(10, 2).Value = "'"
$ws.Cells.Item(10, 2).Style = "Normal"
$ws.Cells.Item(10, 3).Value = "> USD 1,000"
$ws.Cells.Item(10, 4).Value = "'"
$ws.Cells.Item(10, 4).Style = "Normal"

# --- Row 11: Medium ---
$ws.Cells.Item(11, 1).Value = "Medium"
$ws.Cells.Item(11, 2).Value = "<500"
$ws.Cells.Item(11, 3).Value = "< USD 5,000"
$ws.Cells.Item(11, 4).Value = "'"
$ws.Cells.Item(11, 4).Style = "Normal"

# --- Row 12: Large ---
$ws.Cells.Item(12, 1).Value = "Large"
$ws.Cells.Item(12, 2).Value = ">500"
$ws.Cells.Item(12, 3).Value = "> USD 5,000"
$ws.Cells.Item(12, 4).Value = "'"
$ws.Cells.Item(12, 4).Style = "Normal"
